# Recalibration: update Linear and NonLinear sheets with refreshed parameter values.
$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.000112122007561379
$wsLinear.Range("B3").Value = -0.09273814318815624
$wsLinear.Range("B4").Value = 0.001613244288301545
$wsLinear.Range("B5").Value = "[1.0, 0.16014964830042805, 0.006319493385571869, -0.04091551721438996, -0.05858842616238043, -0.013952503734280123, 0.1457823717617883, 0.2997770851256144, 0.12112223282782156, -0.023001484044785407, -0.08167710946294351, -0.06606268930906901, -0.04503916246486379, 0.15204428672204304, 0.30055848564274934, 0.12638316991232157, -0.02583730236550094, -0.09176064262131814, -0.08912174601445778, -0.03584873166383549]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.0003614971820210863
$wsNonLinear.Range("B5").Value = -0.005482692490732029
$wsNonLinear.Range("B6").Value = 0.001676425996407432
$wsNonLinear.Range("B7").Value = 0.004052176248651647
$wsNonLinear.Range("B8").Value = -0.3801356450595673
$wsNonLinear.Range("B9").Value = 0.001544648005416202
$wsNonLinear.Range("B10").Value = "[0.9999999999999998, 0.1633030394619685, 0.009314122979651021, -0.03553475356011592, -0.05403332269742089, -0.016428146199346197, 0.1429110568537933, 0.2966034362638458, 0.11789676010786454, -0.023266859829068506, -0.07719713704763413, -0.06262277265665875, -0.045655405754486816, 0.15025443963985025, 0.2966896631727092, 0.1245880652839908, -0.024112123514402145, -0.08787771097426056, -0.08409609855327078, -0.03490681841406111]"
